$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 164, shifting existing rows 164-179 down to 166-181.
$ws.Range("A164:A165").EntireRow.Insert()

# Fill the two newly inserted rows (164 and 165) with the boilerplate values
# that are constant across every data row in this sheet (A,B,C,E,F,G,H,N,O,Q,R),
# copied from row 166 (the row that used to be row 164 before the insert).
# Row 164
$ws.Cells.Item(164, 1).Value2  = $ws.Cells.Item(166, 1).Value2   # A Mercado ID
$ws.Cells.Item(164, 2).Value2  = $ws.Cells.Item(166, 2).Value2   # B Mercado
$ws.Cells.Item(164, 3).Value2  = $ws.Cells.Item(166, 3).Value2   # C Region
$ws.Cells.Item(164, 5).Value2  = $ws.Cells.Item(166, 5).Value2   # E Codreg
$ws.Cells.Item(164, 6).Value2  = $ws.Cells.Item(166, 6).Value2   # F Categoria ID
$ws.Cells.Item(164, 7).Value2  = $ws.Cells.Item(166, 7).Value2   # G Categoria
$ws.Cells.Item(164, 8).Value2  = $ws.Cells.Item(166, 8).Value2   # H Variedad
$ws.Cells.Item(164, 14).Value2 = $ws.Cells.Item(166, 14).Value2  # N Unidad de comercializacion
$ws.Cells.Item(164, 15).Value2 = $ws.Cells.Item(166, 15).Value2  # O Origen
$ws.Cells.Item(164, 17).Value2 = $ws.Cells.Item(166, 17).Value2  # Q Kg o Unidades
$ws.Cells.Item(164, 18).Value2 = $ws.Cells.Item(166, 18).Value2  # R Clasificacion
$ws.Cells.Item(164, 4).NumberFormat = $ws.Cells.Item(166, 4).NumberFormat  # D date format/style

# Row 165
$ws.Cells.Item(165, 1).Value2  = $ws.Cells.Item(166, 1).Value2
$ws.Cells.Item(165, 2).Value2  = $ws.Cells.Item(166, 2).Value2
$ws.Cells.Item(165, 3).Value2  = $ws.Cells.Item(166, 3).Value2
$ws.Cells.Item(165, 5).Value2  = $ws.Cells.Item(166, 5).Value2
$ws.Cells.Item(165, 6).Value2  = $ws.Cells.Item(166, 6).Value2
$ws.Cells.Item(165, 7).Value2  = $ws.Cells.Item(166, 7).Value2
$ws.Cells.Item(165, 8).Value2  = $ws.Cells.Item(166, 8).Value2
$ws.Cells.Item(165, 14).Value2 = $ws.Cells.Item(166, 14).Value2
$ws.Cells.Item(165, 15).Value2 = $ws.Cells.Item(166, 15).Value2
$ws.Cells.Item(165, 17).Value2 = $ws.Cells.Item(166, 17).Value2
$ws.Cells.Item(165, 18).Value2 = $ws.Cells.Item(166, 18).Value2
$ws.Cells.Item(165, 4).NumberFormat = $ws.Cells.Item(166, 4).NumberFormat

# Row 164 specific values
$ws.Range("D164").Value = 44984
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 300
$ws.Range("K164").Value = 800
$ws.Range("L164").Value = 800
$ws.Range("M164").Value = 800
$ws.Range("P164").Value = 800

# Row 165 specific values
$ws.Range("D165").Value = 44984
$ws.Range("I165").Value = "Segunda"
$ws.Range("J165").Value = 200
$ws.Range("K165").Value = 600
$ws.Range("L165").Value = 600
$ws.Range("M165").Value = 600
$ws.Range("P165").Value = 600
